$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-9 from 45208 (2023-10-09)
# to 45212 (2023-10-13), keeping the existing date format/style intact.
$ws.Range("C2:C9").Value = 45212
